$d = $word.ActiveDocument

$pairs = @(
    @("762÷7=108, 6", "984÷8=123, 0"),
    @("143÷6=23, 5", "848÷8=106, 0"),
    @("145÷7=20, 5", "580÷6=96, 4"),
    @("479÷5=95, 4", "308÷7=44, 0"),
    @("616÷8=77, 0", "829÷7=118, 3"),
    @("589÷3=196, 1", "223÷2=111, 1"),
    @("992÷6=165, 2", "573÷3=191, 0"),
    @("941÷7=134, 3", "862÷2=431, 0"),
    @("812÷7=116, 0", "370÷9=41, 1"),
    @("283÷5=56, 3", "182÷2=91, 0"),
    @("561÷5=112, 1", "723÷4=180, 3"),
    @("456÷3=152, 0", "567÷8=70, 7"),
    @("857÷6=142, 5", "250÷4=62, 2"),
    @("176÷3=58, 2", "422÷6=70, 2"),
    @("455÷3=151, 2", "775÷5=155, 0"),
    @("582÷8=72, 6", "537÷9=59, 6"),
    @("477÷4=119, 1", "250÷5=50, 0"),
    @("754÷5=150, 4", "247÷7=35, 2"),
    @("947÷9=105, 2", "568÷7=81, 1"),
    @("177÷8=22, 1", "352÷3=117, 1"),
    @("679÷2=339, 1", "494÷4=123, 2"),
    @("157÷5=31, 2", "518÷9=57, 5"),
    @("138÷9=15, 3", "148÷8=18, 4"),
    @("551÷2=275, 1", "118÷4=29, 2"),
    @("792÷3=264, 0", "200÷9=22, 2")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
